# Update similarity_scores worksheet:
# - Replace file name strings (Assignment4.pdf / Lab3 files) with Resume.docx / Resume.pdf
# - Shrink the comparison matrix from 3x3 (A1:D4) down to 2x2 (A1:C3)
# - Update the similarity values accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New document names
$name1 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/static/content/temp_files\Hangsihak_Sin_Resume.docx"
$name2 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/static/content/temp_files\Hangsihak_Sin_Resume.pdf"

# Remove the third row/column (document 3) entirely, leaving a 2x2 matrix
$ws.Range("D1:D4").Delete() | Out-Null
$ws.Range("A4:D4").Delete() | Out-Null

# Update header labels (row 1) and row labels (column A)
$ws.Range("B1").Value = $name1
$ws.Range("C1").Value = $name2
$ws.Range("A2").Value = $name1
$ws.Range("A3").Value = $name2

# Update the similarity score matrix values
$ws.Range("B2").Value = 0.9999999999999996
$ws.Range("C2").Value = 0.9962553601710249
$ws.Range("B3").Value = 0.9962553601710249
$ws.Range("C3").Value = 0.9999999403953552

$wb.Save()
